# Change the table style on the "PLENARY" recap slide from the deck's
# custom style to the built-in Table_0-replacement style
# {CE369425-43F4-4A6F-9F1E-F6623F676065}.
#
# (PowerPoint UI equivalent: select the table -> Table Design ribbon ->
#  click the new style swatch in the Table Styles gallery.)

$p = $ppt.ActivePresentation

$oldStyleId = "{BE8641CF-A476-41AB-99F6-5D3EE83A9F90}"
$newStyleId = "{CE369425-43F4-4A6F-9F1E-F6623F676065}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
                Write-Host "Slide $si, Shape $shi - table style" $oldStyleId "->" $table.Style
            }
        }
    }
}
